# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 0.0001488876196638067; C = 0.002777888934908601;  D = 0.1575252929769615;  E = 0.496779210170732;  G = 0.6572312797022659 }
    3 = @{ B = 3.230985683306322;     C = 3099.503889238888;      D = 3.900430680208489;   E = 8.660232485948974;  G = 3115.295538088352 }
    4 = @{ B = 3.230985683306322;     C = 1.667794583268128;      D = 3.900430680208489;   E = 0.496779210170732;  G = 9.295990156953671 }
    5 = @{ B = 1.459612070389937;     C = 1.667794583268128;      D = 26.21740644021617;   E = 0.496779210170732;  G = 29.84159230404497 }
    6 = @{ B = 0.6753301551942219;    C = 1.667794583268128;      D = 0.8054896365839992;  E = 0.496779210170732;  G = 3.645393585217082 }
    7 = @{ B = 3.230985683306322;     C = 1.667794583268128;      D = 0.1575252929769615;  E = 0.496779210170732;  G = 5.553084769722144 }
    8 = @{ B = 3.230985683306322;     C = 1.667794583268128;      D = 0.1575252929769615;  E = 0.496779210170732;  G = 5.553084769722144 }
    9 = @{ B = 3.230985683306322;     C = 1.667794583268128;      D = 0.8054896365839992;  E = 0.496779210170732;  G = 6.201049113329182 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
